# The document contains one table with 5 data rows (rows 1, 5, 9, 13, 17 of
# the 20-row table -- the other rows are blank spacer rows) and 5 columns of
# "NNN÷N=" division problems. Each cell is addressed directly by its
# (row, column) position rather than via Find/Replace, because a couple of
# the new values happen to equal other cells' old values (e.g. 434÷7= becomes
# 708÷6=, while the original 708÷6= cell becomes 343÷6=); a global
# find-and-replace run in sequence could re-match text that a previous
# replacement just introduced. Addressing each cell by position avoids that
# collision entirely.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "330÷3="   # was 411÷7=
$t.Cell(1,2).Range.Text = "469÷3="   # was 571÷9=
$t.Cell(1,3).Range.Text = "592÷4="   # was 896÷6=
$t.Cell(1,4).Range.Text = "112÷2="   # was 632÷6=
$t.Cell(1,5).Range.Text = "888÷6="   # was 939÷5=

# Row 5
$t.Cell(5,1).Range.Text = "708÷6="   # was 434÷7=
$t.Cell(5,2).Range.Text = "392÷9="   # was 878÷8=
$t.Cell(5,3).Range.Text = "320÷9="   # was 518÷6=
$t.Cell(5,4).Range.Text = "746÷8="   # was 651÷3=
$t.Cell(5,5).Range.Text = "498÷9="   # was 601÷2=

# Row 9
$t.Cell(9,1).Range.Text = "708÷7="   # was 832÷8=
$t.Cell(9,2).Range.Text = "972÷2="   # was 887÷3=
$t.Cell(9,3).Range.Text = "357÷7="   # was 739÷7=
$t.Cell(9,4).Range.Text = "440÷9="   # was 824÷9=
$t.Cell(9,5).Range.Text = "343÷6="   # was 708÷6=

# Row 13
$t.Cell(13,1).Range.Text = "601÷8="  # was 864÷2=
$t.Cell(13,2).Range.Text = "648÷2="  # was 745÷7=
$t.Cell(13,3).Range.Text = "455÷8="  # was 297÷9=
$t.Cell(13,4).Range.Text = "403÷5="  # was 497÷5=
$t.Cell(13,5).Range.Text = "233÷2="  # was 136÷8=

# Row 17
$t.Cell(17,1).Range.Text = "722÷8="  # was 938÷3=
$t.Cell(17,2).Range.Text = "924÷5="  # was 197÷2=
$t.Cell(17,3).Range.Text = "663÷8="  # was 314÷3=
$t.Cell(17,4).Range.Text = "203÷3="  # was 255÷7=
$t.Cell(17,5).Range.Text = "114÷7="  # was 400÷6=
